# Fill in the "Beat Vegas?" column for the last 5 previously-unlabeled
# games (rows 67-71), then append the results of the games the model was
# run against on Jan 15 2021 (rows 72-79).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G ("Beat Vegas?") for rows 67-71 ---
$ws.Range("G67").Value = "No"
$ws.Range("G68").Value = "Yes"
$ws.Range("G69").Value = "Yes"
$ws.Range("G70").Value = "No"
$ws.Range("G71").Value = "No"

# --- New games for 2021-01-15 (rows 72-79) ---
$newGames = @(
    @("CLE", "NYK",  2,    2.1,  -0.1000000000000001),
    @("BOS", "ORL", -2.5, -8.4,   5.9),
    @("MIL", "DAL", -7,  -13.6,   6.6),
    @("OKC", "CHI",  0,   10.9, -10.9),
    @("MIN", "MEM", -1,   -0.2,  -0.8),
    @("UTA", "ATL", -6,   -9.6,   3.6),
    @("LAL", "NOP", -10, -12.7,   2.699999999999999),
    @("SAC", "LAC",  7,   -5.8,  12.8)
)

$firstRow = 72
$lastRow = $firstRow + $newGames.Count - 1

# Write the date serial number (midnight 2021-01-15) first, then clone the
# date format from an existing date cell so the new cells share the same
# "yyyy-mm-dd" style (s="2") instead of getting a brand-new number format.
$ws.Range("A$($firstRow):A$($lastRow)").Value = 44211
$ws.Range("A71").Copy()
$ws.Range("A$($firstRow):A$($lastRow)").PasteSpecial(-4122)

$row = $firstRow
foreach ($game in $newGames) {
    $ws.Range("B$row").Value = $game[0]
    $ws.Range("C$row").Value = $game[1]
    $ws.Range("D$row").Value = $game[2]
    $ws.Range("E$row").Value = $game[3]
    $ws.Range("F$row").Value = $game[4]
    $row++
}
